# "Modificacion y organización de muestreo"
#
# The external workbook ("Muestreo de datos victus residencias nuevo .xlsx")
# that this file links to (external link [1]) had its "Inmueble" and
# "Residente" sheets removed/reorganized. As a result, the cached formulas in
# the "Residente" sheet here that pulled the pre-computed "inmueble" label
# from that external workbook (via [1]Residente!J4, J5, J6) are now broken
# references. Re-create that breakage here: those three formulas now
# evaluate to #REF! errors, which ripples into the dependent K4:K6
# "Combinación única" CONCAT formulas as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Residente")

$ws.Range("J4").Formula = "=#REF!"
$ws.Range("J5").Formula = "=#REF!"
$ws.Range("J6").Formula = "=#REF!"
